$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("D29").Value = 400
$ws.Range("D30").Value = 400
$ws.Range("D31").Value = 400
$ws.Range("D32").Value = 400
$ws.Range("D33").Value = 1182
